# Auto-generated script applying scheduled market-data refresh to Jenova_Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per refreshed Universalis data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 640.5
$ws.Range("I103").Value = 864.6667
$ws.Range("J103").Value = 416.33334
$ws.Range("K103").Value = 2594.0001
$ws.Range("L103").Value = 1249.00002
$ws.Range("M103").Value = -2008.0001
$ws.Range("N103").Value = -2421.00002
$ws.Range("H113").Value = 6000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -12508
$ws.Range("H132").Value = 3014.4849
$ws.Range("I132").Value = 2785
$ws.Range("J132").Value = 4299.6
$ws.Range("K132").Value = 8355
$ws.Range("L132").Value = 12898.8
$ws.Range("M132").Value = -5825
$ws.Range("N132").Value = -17958.8
$ws.Range("H137").Value = 716909.9399999999
$ws.Range("J137").Value = 1004434.9
$ws.Range("L137").Value = 3013304.7
$ws.Range("N137").Value = -3018404.7
$ws.Range("H138").Value = 6908.8774
$ws.Range("J138").Value = 7902.125
$ws.Range("L138").Value = 23706.375
$ws.Range("N138").Value = -33986.375
$ws.Range("H141").Value = 2953.5454
$ws.Range("I141").Value = 2953.5454
$ws.Range("K141").Value = 8860.636200000001
$ws.Range("M141").Value = -3680.636200000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3595.6714
$ws.Range("I32").Value = 3595.6714
$ws.Range("K32").Value = 3595.6714
$ws.Range("M32").Value = -3308.6714
$ws.Range("H45").Value = 2722.6
$ws.Range("I45").Value = 1905.1
$ws.Range("J45").Value = 4357.6
$ws.Range("K45").Value = 1905.1
$ws.Range("L45").Value = 4357.6
$ws.Range("M45").Value = -1528.1
$ws.Range("N45").Value = -5111.6
$ws.Range("H62").Value = 26331.666
$ws.Range("I62").Value = 29000
$ws.Range("J62").Value = 24997.5
$ws.Range("K62").Value = 29000
$ws.Range("L62").Value = 24997.5
$ws.Range("M62").Value = -28376
$ws.Range("N62").Value = -26245.5
$ws.Range("H65").Value = 26331.666
$ws.Range("I65").Value = 29000
$ws.Range("J65").Value = 24997.5
$ws.Range("K65").Value = 87000
$ws.Range("L65").Value = 74992.5
$ws.Range("M65").Value = -83880
$ws.Range("N65").Value = -81232.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1158
$ws.Range("J64").Value = 863.3333
$ws.Range("L64").Value = 863.3333
$ws.Range("N64").Value = -1313.3333
$ws.Range("H67").Value = 1158
$ws.Range("J67").Value = 863.3333
$ws.Range("L67").Value = 863.3333
$ws.Range("N67").Value = -2423.3333
$ws.Range("H94").Value = 1080.3846
$ws.Range("I94").Value = 1054.5
$ws.Range("J94").Value = 1166.6666
$ws.Range("K94").Value = 1054.5
$ws.Range("L94").Value = 1166.6666
$ws.Range("M94").Value = -603.5
$ws.Range("N94").Value = -2068.6666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1481.25
$ws.Range("I16").Value = 1469
$ws.Range("J16").Value = 1493.5
$ws.Range("K16").Value = 1469
$ws.Range("L16").Value = 1493.5
$ws.Range("M16").Value = -1182
$ws.Range("N16").Value = -2067.5
$ws.Range("H31").Value = 24932.137
$ws.Range("I31").Value = 1594.2424
$ws.Range("J31").Value = 94945.82000000001
$ws.Range("K31").Value = 1594.2424
$ws.Range("L31").Value = 94945.82000000001
$ws.Range("M31").Value = -1299.2424
$ws.Range("N31").Value = -95535.82000000001
$ws.Range("H34").Value = 24932.137
$ws.Range("I34").Value = 1594.2424
$ws.Range("J34").Value = 94945.82000000001
$ws.Range("K34").Value = 1594.2424
$ws.Range("L34").Value = 94945.82000000001
$ws.Range("M34").Value = -1392.2424
$ws.Range("N34").Value = -95349.82000000001
$ws.Range("H62").Value = 4178.8
$ws.Range("J62").Value = 5298.3335
$ws.Range("L62").Value = 5298.3335
$ws.Range("N62").Value = -6546.3335
$ws.Range("H65").Value = 4178.8
$ws.Range("J65").Value = 5298.3335
$ws.Range("L65").Value = 26491.6675
$ws.Range("N65").Value = -32731.6675
$ws.Range("H99").Value = 6482.4375
$ws.Range("I99").Value = 5120.5
$ws.Range("K99").Value = 5120.5
$ws.Range("M99").Value = -3622.5
$ws.Range("H113").Value = 1481.25
$ws.Range("I113").Value = 1469
$ws.Range("J113").Value = 1493.5
$ws.Range("K113").Value = 1469
$ws.Range("L113").Value = 1493.5
$ws.Range("M113").Value = 701
$ws.Range("N113").Value = -5833.5
$ws.Range("H126").Value = 6482.4375
$ws.Range("I126").Value = 5120.5
$ws.Range("K126").Value = 15361.5
$ws.Range("M126").Value = -12891.5
$ws.Range("H132").Value = 4547.793
$ws.Range("I132").Value = 4112.0527
$ws.Range("K132").Value = 12336.1581
$ws.Range("M132").Value = -9806.158100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 31643.115
$ws.Range("I5").Value = 50324.75
$ws.Range("K5").Value = 150974.25
$ws.Range("M5").Value = -150862.25
$ws.Range("H56").Value = 6499.75
$ws.Range("I56").Value = 6499.75
$ws.Range("K56").Value = 6499.75
$ws.Range("M56").Value = -5969.75
$ws.Range("H62").Value = 14664.833
$ws.Range("J62").Value = 14997.8
$ws.Range("L62").Value = 44993.39999999999
$ws.Range("N62").Value = -46365.39999999999
$ws.Range("H65").Value = 14664.833
$ws.Range("J65").Value = 14997.8
$ws.Range("L65").Value = 134980.2
$ws.Range("N65").Value = -141844.2
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 10015
$ws.Range("J74").Value = 10015
$ws.Range("L74").Value = 30045
$ws.Range("N74").Value = -32167
$ws.Range("H75").Value = 1494.6666
$ws.Range("I75").Value = 1494.6666
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 4483.9998
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -3485.9998
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 10015
$ws.Range("J77").Value = 10015
$ws.Range("L77").Value = 90135
$ws.Range("N77").Value = -100743
$ws.Range("H78").Value = 1494.6666
$ws.Range("I78").Value = 1494.6666
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 13451.9994
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -8459.999400000001
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 502
$ws.Range("I80").Value = 502
$ws.Range("K80").Value = 1506
$ws.Range("M80").Value = -570
$ws.Range("H81").Value = 206723.39
$ws.Range("J81").Value = 206723.39
$ws.Range("L81").Value = 620170.17
$ws.Range("N81").Value = -622416.17
$ws.Range("H83").Value = 502
$ws.Range("I83").Value = 502
$ws.Range("K83").Value = 4518
$ws.Range("M83").Value = 162
$ws.Range("H84").Value = 206723.39
$ws.Range("J84").Value = 206723.39
$ws.Range("L84").Value = 1860510.51
$ws.Range("N84").Value = -1871742.51
$ws.Range("H88").Value = 7999
$ws.Range("J88").Value = 7999
$ws.Range("L88").Value = 23997
$ws.Range("N88").Value = -24853
$ws.Range("H91").Value = 7999
$ws.Range("J91").Value = 7999
$ws.Range("L91").Value = 23997
$ws.Range("N91").Value = -26961
$ws.Range("H122").Value = 50824.75
$ws.Range("I122").Value = 723.8570999999999
$ws.Range("K122").Value = 6514.7139
$ws.Range("M122").Value = -4064.7139
$ws.Range("H135").Value = 31643.115
$ws.Range("I135").Value = 50324.75
$ws.Range("K135").Value = 452922.75
$ws.Range("M135").Value = -450387.75
$ws.Range("H136").Value = 5430.625
$ws.Range("I136").Value = 2889.4
$ws.Range("J136").Value = 9666
$ws.Range("K136").Value = 8668.200000000001
$ws.Range("L136").Value = 28998
$ws.Range("M136").Value = -3568.200000000001
$ws.Range("N136").Value = -39198
$ws.Range("H139").Value = 5427.2383
$ws.Range("I139").Value = 4625.0625
$ws.Range("K139").Value = 13875.1875
$ws.Range("M139").Value = -8735.1875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H103").Value = 33333
$ws.Range("J103").Value = 33333
$ws.Range("L103").Value = 33333
$ws.Range("N103").Value = -35677
$ws.Range("H122").Value = 3185.0715
$ws.Range("I122").Value = 1410.1538
$ws.Range("J122").Value = 4723.3335
$ws.Range("K122").Value = 4230.4614
$ws.Range("L122").Value = 14170.0005
$ws.Range("M122").Value = -1780.4614
$ws.Range("N122").Value = -19070.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5154.364
$ws.Range("J46").Value = 3955.2222
$ws.Range("L46").Value = 3955.2222
$ws.Range("N46").Value = -4331.2222

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7710.222
$ws.Range("I41").Value = 2498
$ws.Range("J41").Value = 8361.75
$ws.Range("K41").Value = 2498
$ws.Range("L41").Value = 8361.75
$ws.Range("M41").Value = -2108
$ws.Range("N41").Value = -9141.75

Write-Host "Applied market-data refresh across all profession sheets."